# Reverse words in a string
#
# Adds a new LeetCode entry ("Reverse Words in a String") plus two more rows
# ("Two Sum" and "Rotate Array") to the tracker sheet, and extends the tags
# on the first existing row.
#
# xlPasteFormats = -4122 (used so pasted cells reuse an *existing* cellXf
# instead of Hyperlinks.Add's own freshly-minted one).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$xlPasteFormats = -4122

# --- Make room: insert two blank rows right after row 3. ---
# This shifts the old row 4 ("Array" / print-words-vertically, untouched)
# down to row 6, and old rows 5..13 down to rows 7..15.
$ws.Rows("4:5").Insert() | Out-Null

# Excel auto-fills the new rows' formatting from the row above (row 3), so
# B4/D4 and B5/D5 already picked up styles 2/3. Row 5 has no D column entry
# in the final layout, so drop that leftover formatting entirely.
$ws.Range("D5").Clear() | Out-Null

# --- Row 7 (previously row 5): Two Sum ---
$ws.Hyperlinks.Add($ws.Range("B7"), "https://leetcode.com/problems/two-sum/") | Out-Null
$ws.Range("B3").Copy() | Out-Null
$ws.Range("B7").PasteSpecial($xlPasteFormats) | Out-Null
$ws.Range("A7").Value = "Array, Hashmap"
$ws.Range("C7").Value = "Solved using BruteForce. Can be improved"
$ws.Range("D3").Copy() | Out-Null
$ws.Range("D7").PasteSpecial($xlPasteFormats) | Out-Null

# --- Row 4 (new): Reverse Words in a String ---
$ws.Hyperlinks.Add($ws.Range("B4"), "https://leetcode.com/problems/reverse-words-in-a-string/") | Out-Null
$ws.Range("B3").Copy() | Out-Null
$ws.Range("B4").PasteSpecial($xlPasteFormats) | Out-Null
$ws.Range("C4").Value = "Tokenize the string and reverse"
$ws.Range("A4").Value = "Array, String"

# --- Row 3: extend the existing tag list ---
$ws.Range("A3").Value = "Array, Two Pointer, Recursion, String"

# --- Row 5 (new): Rotate Array ---
$ws.Hyperlinks.Add($ws.Range("B5"), "https://leetcode.com/problems/rotate-array/") | Out-Null
$ws.Range("B3").Copy() | Out-Null
$ws.Range("B5").PasteSpecial($xlPasteFormats) | Out-Null
$ws.Range("A5").Value = "Array"

# --- Two extra trailing "Array" rows ---
$ws.Range("A14").Value = "Array"
$ws.Range("A15").Value = "Array"

# --- Match the author's last selection ---
$ws.Range("B17").Select() | Out-Null

Write-Host "edit applied"
